$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

function Set-NumberCell($addr, $val) {
    $ws.Range($addr).Value = $val
}

function Set-EmptyCell($addr) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = ""
    $rng.Style = "Normal"
}

# --- Update header row (row 1) ---
$ws.Range("AE1").Value = "Today  Date"
$ws.Range("AF1").Value = "Today Date"
$ws.Range("AE1").Copy()
$ws.Range("AF1").PasteSpecial(-4122)

# --- Row 2 data: DHARMANA DHANA LAKSHMI ---
Set-TextCell   "A2"  "DHARMANA DHANA LAKSHMI"
Set-TextCell   "B2"  "Dhana"
Set-TextCell   "C2"  "Laxman"
Set-TextCell   "D2"  "Tulasi"
Set-NumberCell "E2"  9493541829
Set-TextCell   "F2"  "dharmanadhanalaxmi006@outlook.com"
Set-TextCell   "G2"  "Andhra Pradesh state, Srikakulam district."
Set-NumberCell "H2"  463426316719
Set-TextCell   "I2"  "Father"
Set-NumberCell "J2"  9550791829
Set-TextCell   "K2"  "Intern"
Set-TextCell   "L2"  "Data Science"
Set-TextCell   "M2"  "2025-04-16"
Set-TextCell   "N2"  "Full-Time"
Set-TextCell   "O2"  "Ganesh"
Set-TextCell   "P2"  "Remote"
Set-NumberCell "Q2"  1253647890
Set-TextCell   "R2"  "Dhana Lakshmi"
Set-TextCell   "S2"  "XXXXXXXXX10"
Set-TextCell   "T2"  "Laptop"
Set-TextCell   "U2"  "Canva"
Set-TextCell   "V2"  "Headset"
Set-TextCell   "W2"  "Yes"
Set-TextCell   "X2"  "Yes"
Set-TextCell   "Y2"  "Yes"
Set-TextCell   "Z2"  "Yes"
Set-TextCell   "AA2" "Yes"
Set-TextCell   "AB2" "Yes"
Set-EmptyCell  "AC2"
Set-TextCell   "AD2" "Dhana Lakshmi"
Set-EmptyCell  "AE2"
Set-EmptyCell  "AF2"

# --- Row 3 data: Tejaswnini ---
Set-TextCell  "A3"  "Tejaswnini"
Set-TextCell  "B3"  "Teju"
Set-TextCell  "C3"  "Ramana"
Set-TextCell  "D3"  "Sita"
Set-TextCell  "E3"  "7675993724"
Set-TextCell  "F3"  "tejaswanipulugu171@gmail.com"
Set-TextCell  "G3"  "amalapuram"
Set-TextCell  "H3"  "123456789123"
Set-TextCell  "I3"  "Dhana"
Set-TextCell  "J3"  "9493541829"
Set-TextCell  "K3"  "Intern"
Set-TextCell  "L3"  "Data Science"
Set-TextCell  "M3"  "2025-04-16"
Set-TextCell  "N3"  "Full-Time"
Set-TextCell  "O3"  "John Kiran"
Set-TextCell  "P3"  "Office"
Set-TextCell  "Q3"  "1253647890"
Set-TextCell  "R3"  "Dhana Lakshmi"
Set-TextCell  "S3"  "XXXXXXXXX10"
Set-TextCell  "T3"  "Laptop"
Set-TextCell  "U3"  "Canva"
Set-TextCell  "V3"  "Headset"
Set-TextCell  "W3"  "Yes"
Set-TextCell  "X3"  "Yes"
Set-TextCell  "Y3"  "Yes"
Set-TextCell  "Z3"  "Yes"
Set-TextCell  "AA3" "Yes"
Set-TextCell  "AB3" "Yes"
Set-TextCell  "AC3" "NA"
Set-TextCell  "AD3" "Tejaswini"
Set-EmptyCell "AE3"
Set-EmptyCell "AF3"
